$d = $word.ActiveDocument

# 1. Update the timestamp in the date line: 09:01:00 PM -> 09:01:33 PM.
#    Locate the paragraph using the "Date" style rather than a hard-coded
#    index so the edit is resilient to unrelated structural changes.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $p.Range.Text = "June  16, 2021 (09:01:33 PM)"
        break
    }
}

# 2. In the first table, the header-ish row reads:
#    Input | returns | value | returns value | (empty)
#    Split the "returns value" cell into "returns" (kept in place) and a
#    new "value" cell (centered) in the following column.
$table = $d.Tables.Item(1)

$returnsCell = $null
for ($r = 1; $r -le $table.Rows.Count; $r++) {
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $cell = $table.Cell($r, $c)
        $text = $cell.Range.Text.TrimEnd([char]7, [char]13)
        if ($text -eq "returns value") {
            $returnsCell = $cell
            $returnsRow = $r
            $returnsCol = $c
        }
    }
}

if ($returnsCell -ne $null) {
    $returnsCell.Range.Text = "returns"

    $valueCell = $table.Cell($returnsRow, $returnsCol + 1)
    $valueCell.Range.ParagraphFormat.Alignment = 1
    $valueCell.Range.Text = "value"
}
